# Updated cryptos list on Mon Jul 31 10:57:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a TEXT string
# (many "price" values look numeric, e.g. "0.7047" or "1.000", and Excel's
# normal type inference on Range.Value would silently coerce them to a
# number and drop the trailing zeros / formatting). Flip the cell to text
# format just long enough to type the value in as a string, then clear the
# formatting back off so the cell keeps the default (unstyled) look.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.400.20"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.869.20"
$ws.Range("E3").Value = "  -0.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
Set-TextValue "D5" "243.63"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6 - XRP
Set-TextValue "D6" "0.7047"
$ws.Range("E6").Value = "  -2.47%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.07922"
$ws.Range("E8").Value = "  -1.01%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3133"
$ws.Range("E9").Value = "  -0.56%  "

# Row 10 - Solana
Set-TextValue "D10" "24.52"
$ws.Range("E10").Value = "  -1.57%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07841"
$ws.Range("E11").Value = "  -4.55%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.898.31"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13 - Litecoin
Set-TextValue "D13" "93.79"
$ws.Range("E13").Value = "  -0.65%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.170"
$ws.Range("E14").Value = "  -0.93%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.7011"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16 - Uniswap
Set-TextValue "D16" "6.516"
$ws.Range("E16").Value = "  +2.10%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.000008401"
$ws.Range("E17").Value = "  -0.87%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.481.53"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "252.52"
$ws.Range("E19").Value = "  +3.90%  "

# Row 20 - WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "2.147.60"
$ws.Range("E20").Value = "  +0.48%  "

# Row 22 - Dai
Set-TextValue "D22" "1.001"
$ws.Range("E22").Value = "  -0.07%  "

# Row 23 - Chainlink
Set-TextValue "D23" "7.677"
$ws.Range("E23").Value = "  -1.11%  "

# Row 24 - BinanceUSD
Set-TextValue "D24" "1.000"
$ws.Range("E24").Value = "  -0.23%  "

# Row 25 - Stellar
Set-TextValue "D25" "0.1552"
$ws.Range("E25").Value = "  -3.21%  "

# Row 27 - Monero
Set-TextValue "D27" "161.71"
$ws.Range("E27").Value = "  -0.51%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +1.77%  "

# Row 29 - PancakeSwap
Set-TextValue "D29" "1.506"
$ws.Range("E29").Value = "  +0.30%  "

# Row 30 - Filecoin
Set-TextValue "D30" "4.316"
$ws.Range("E30").Value = "  -1.98%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "4.260"
$ws.Range("E31").Value = "  -0.92%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  +2.39%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.05271"
$ws.Range("E33").Value = "  -1.56%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.895"
$ws.Range("E34").Value = "  -1.98%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.179"
$ws.Range("E35").Value = "  +0.35%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.7507"
$ws.Range("E36").Value = "  -0.99%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.712"
$ws.Range("E37").Value = "  +0.47%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.06%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.275.27"
$ws.Range("E39").Value = "  -0.61%  "

# Row 40 - MXToken
Set-TextValue "D40" "2.771"
$ws.Range("E40").Value = "  +0.63%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "0.8915"
$ws.Range("E41").Value = "  -1.84%  "

# Row 42 - was Quant, now FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "6.039"
$ws.Range("E42").Value = "  -6.08%  "

# Row 43 - was FraxShare, now Quant
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "109.17"
$ws.Range("E43").Value = "  -3.22%  "

# Row 44 - Aave
Set-TextValue "D44" "70.88"
$ws.Range("E44").Value = "  -4.40%  "

# Row 45 - PaxDollar
Set-TextValue "D45" "1.001"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "2.044.30"
$ws.Range("E46").Value = "  +0.57%  "

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  -4.59%  "

# Row 48 - RenderToken
Set-TextValue "D48" "1.804"
$ws.Range("E48").Value = "  +0.61%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "9.615"
$ws.Range("E49").Value = "  +1.52%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.89%  "

# Row 51 - TheSandbox
Set-TextValue "D51" "0.4299"
$ws.Range("E51").Value = "  -0.83%  "
